$d = $word.ActiveDocument

# Hunk 1: insert new Heading2 "Small SW safety measure" section + body paragraphs before "New requirement" heading
$rng = $d.Content
$rng.Find.Execute("New requirement", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$p = $rng.Paragraphs(1)
$insertPoint = $d.Range($p.Range.Start, $p.Range.Start)
$xml1 = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="Heading2"/><w:rPr><w:lang w:val="en-US" w:bidi="ar-JO"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US" w:bidi="ar-JO"/></w:rPr><w:t>Small SW safety measure – 21/11/22</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:lang w:val="en-US" w:bidi="ar-JO"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US" w:bidi="ar-JO"/></w:rPr><w:t xml:space="preserve">Prepare Wage Table every time before rights </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:lang w:val="en-US" w:bidi="ar-JO"/></w:rPr><w:t>are</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:lang w:val="en-US" w:bidi="ar-JO"/></w:rPr><w:t xml:space="preserve"> computed.</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:lang w:val="en-US" w:bidi="ar-JO"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US" w:bidi="ar-JO"/></w:rPr><w:t>Just to be sure there is no way to use a table that is not up-to-date.</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:lang w:val="en-US" w:bidi="ar-JO"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US" w:bidi="ar-JO"/></w:rPr><w:t>Version changed to “</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New"/><w:lang w:val="en-US"/></w:rPr><w:t>v1.2.</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New"/><w:lang w:val="en-US"/></w:rPr><w:t>1.beta</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New"/><w:lang w:val="en-US"/></w:rPr><w:t>2</w:t></w:r><w:r><w:rPr><w:lang w:val="en-US" w:bidi="ar-JO"/></w:rPr><w:t>”</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:lang w:val="en-US" w:bidi="ar-JO"/></w:rPr></w:pPr></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$insertPoint.InsertXML($xml1)

# Hunk 2: split the older "v1.2.1.beta1" text run into 3 runs (text unchanged) with proofErr markers
$rng2 = $d.Content
$rng2.Find.Execute("v1.2.1.beta1", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$p2 = $rng2.Paragraphs(1)
$pr2 = $p2.Range
$contentRange2 = $d.Range($pr2.Start, $pr2.End - 1)
$xml2 = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:lang w:val="en-US" w:bidi="ar-JO"/></w:rPr><w:t>Version changed to “</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New"/><w:lang w:val="en-US"/></w:rPr><w:t>v1.2.</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New"/><w:lang w:val="en-US"/></w:rPr><w:t>1.beta</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New"/><w:lang w:val="en-US"/></w:rPr><w:t>1</w:t></w:r><w:r><w:rPr><w:lang w:val="en-US" w:bidi="ar-JO"/></w:rPr><w:t>”</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$contentRange2.InsertXML($xml2)

# Hunk 3: add lastRenderedPageBreak to "Save and restore all types of wage periods"
$rng3 = $d.Content
$rng3.Find.Execute("Save and restore all types of wage periods", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$p3 = $rng3.Paragraphs(1)
$pr3 = $p3.Range
$contentRange3 = $d.Range($pr3.Start, $pr3.End - 1)
$xml3 = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:lastRenderedPageBreak/><w:t>Save and restore all types of wage periods</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$contentRange3.InsertXML($xml3)

# Hunk 4: remove lastRenderedPageBreak from "Truthfulness of all wage-related computations"
$rng4 = $d.Content
$rng4.Find.Execute("Truthfulness of all wage-related computations", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$p4 = $rng4.Paragraphs(1)
$pr4 = $p4.Range
$contentRange4 = $d.Range($pr4.Start, $pr4.End - 1)
$xml4 = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Truthfulness of all wage-related computations</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$contentRange4.InsertXML($xml4)

# Hunk 5: add lastRenderedPageBreak to "Handle single period definition from work-period dialog"
$rng5 = $d.Content
$rng5.Find.Execute("Handle single period definition from work-period dialog", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$p5 = $rng5.Paragraphs(1)
$pr5 = $p5.Range
$contentRange5 = $d.Range($pr5.Start, $pr5.End - 1)
$xml5 = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:lastRenderedPageBreak/><w:t>Handle single period definition from work-period dialog</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$contentRange5.InsertXML($xml5)

# Hunk 6: remove lastRenderedPageBreak from "All holiday files updated with data for 2022"
$rng6 = $d.Content
$rng6.Find.Execute("All holiday files updated with data for 2022", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$p6 = $rng6.Paragraphs(1)
$pr6 = $p6.Range
$contentRange6 = $d.Range($pr6.Start, $pr6.End - 1)
$xml6 = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>All holiday files updated with data for 2022 – 26/9/22</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$contentRange6.InsertXML($xml6)

Write-Host "All hunks applied"
